# Update handout and presentation (time data and video embedded)
#
# 1) Team-members time table ("Tabelle 13"): fill in the TODO hour values.
# 2) Task-type time table ("Tabelle 15"): fill in the TODO hour values.
# 3) Refresh the cached "datetimeFigureOut" date placeholders (slide master,
#    every slide layout, and the notes master) from 25.06.2020 -> 28.06.2020.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------------
# 1) Team-members table (id=13 "Tabelle 13") - shape index 6 on the slide.
# ---------------------------------------------------------------------------
$membersTable = $s.Shapes.Item(6).Table

# Row 2 = Fabian : Overall | Project Mgmt | Analysis & Design | Implementation | Testing
$membersTable.Cell(2, 2).Shape.TextFrame.TextRange.Text = "156,5"
$membersTable.Cell(2, 3).Shape.TextFrame.TextRange.Text = "19,25"
$membersTable.Cell(2, 4).Shape.TextFrame.TextRange.Text = "20,25"
$membersTable.Cell(2, 5).Shape.TextFrame.TextRange.Text = "109"
$membersTable.Cell(2, 6).Shape.TextFrame.TextRange.Text = "8"

# Row 3 = Luca
$membersTable.Cell(3, 2).Shape.TextFrame.TextRange.Text = "106,25"
$membersTable.Cell(3, 3).Shape.TextFrame.TextRange.Text = "26,25"
$membersTable.Cell(3, 4).Shape.TextFrame.TextRange.Text = "27,25"
$membersTable.Cell(3, 5).Shape.TextFrame.TextRange.Text = "47,25"
$membersTable.Cell(3, 6).Shape.TextFrame.TextRange.Text = "5,5"

# Row 4 = Nicolas
$membersTable.Cell(4, 2).Shape.TextFrame.TextRange.Text = "194,16"
$membersTable.Cell(4, 3).Shape.TextFrame.TextRange.Text = "67,91"
$membersTable.Cell(4, 4).Shape.TextFrame.TextRange.Text = "42,75"
$membersTable.Cell(4, 5).Shape.TextFrame.TextRange.Text = "56,75"
$membersTable.Cell(4, 6).Shape.TextFrame.TextRange.Text = "26,75"

# ---------------------------------------------------------------------------
# 2) Task-type table (id=15 "Tabelle 15") - shape index 7 on the slide.
# ---------------------------------------------------------------------------
$taskTable = $s.Shapes.Item(7).Table

$taskTable.Cell(2, 2).Shape.TextFrame.TextRange.Text = "113,41"  # Project Management
$taskTable.Cell(3, 2).Shape.TextFrame.TextRange.Text = "90,25"   # Analysis & Design
$taskTable.Cell(4, 2).Shape.TextFrame.TextRange.Text = "213"     # Implementation
$taskTable.Cell(5, 2).Shape.TextFrame.TextRange.Text = "40,25"   # Testing

# ---------------------------------------------------------------------------
# 3) Refresh the cached date placeholders: 25.06.2020 -> 28.06.2020
# ---------------------------------------------------------------------------
$newDate = "28.06.2020"

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq "25.06.2020") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

$layouts = $master.CustomLayouts
for ($l = 1; $l -le $layouts.Count; $l++) {
    $layout = $layouts.Item($l)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq "25.06.2020") {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

$notesMaster = $p.NotesMaster
for ($i = 1; $i -le $notesMaster.Shapes.Count; $i++) {
    $sh = $notesMaster.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq "25.06.2020") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}
